$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.2448979591836735
$ws.Range("C2").Value2 = 0.4897959183673469
$ws.Range("J2").Value2 = 0.08163265306122448
$ws.Range("P2").Value2 = 0.1020408163265306
$ws.Range("S2").Value2 = 0.08163265306122448
$ws.Range("J3").Value2 = 0.16
$ws.Range("P3").Value2 = 0.6
$ws.Range("S3").Value2 = 0.24
$ws.Range("J4").Value2 = 0.1428571428571428
$ws.Range("P4").Value2 = 0.8571428571428571
$ws.Range("B6").Value2 = 0.05
$ws.Range("F6").Value2 = 0.05
$ws.Range("J6").Value2 = 0.4
$ws.Range("O6").Value2 = 0.05
$ws.Range("Q6").Value2 = 0.1
$ws.Range("R6").Value2 = 0.05
$ws.Range("S6").Value2 = 0.3
$ws.Range("B7").Value2 = 0.05555555555555555
$ws.Range("F7").Value2 = 0.05555555555555555
$ws.Range("J7").Value2 = 0.1666666666666667
$ws.Range("O7").Value2 = 0.05555555555555555
$ws.Range("Q7").Value2 = 0.2222222222222222
$ws.Range("S7").Value2 = 0.2777777777777778
$ws.Range("B8").Value2 = 0.04651162790697674
$ws.Range("D8").Value2 = 0.02325581395348837
$ws.Range("J8").Value2 = 0.1627906976744186
$ws.Range("O8").Value2 = 0.02325581395348837
$ws.Range("Q8").Value2 = 0.09302325581395349
$ws.Range("R8").Value2 = 0.2790697674418605
$ws.Range("S8").Value2 = 0.3720930232558139
$ws.Range("B9").Value2 = 0.09090909090909091
$ws.Range("J9").Value2 = 0.1818181818181818
$ws.Range("R9").Value2 = 0.1818181818181818
$ws.Range("S9").Value2 = 0.5454545454545454
$ws.Range("B10").Value2 = 0.1042345276872964
$ws.Range("D10").Value2 = 0.01954397394136808
$ws.Range("F10").Value2 = 0.02931596091205212
$ws.Range("J10").Value2 = 0.1758957654723127
$ws.Range("O10").Value2 = 0.006514657980456026
$ws.Range("Q10").Value2 = 0.247557003257329
$ws.Range("R10").Value2 = 0.1335504885993485
$ws.Range("S10").Value2 = 0.2833876221498371
$ws.Range("G11").Value2 = 0.2068965517241379
$ws.Range("J11").Value2 = 0.03448275862068965
$ws.Range("K11").Value2 = 0.2068965517241379
$ws.Range("L11").Value2 = 0.5517241379310345
$ws.Range("G12").Value2 = 0.6875
$ws.Range("J12").Value2 = 0.25
$ws.Range("S12").Value2 = 0.0625
$ws.Range("G13").Value2 = 0.6666666666666666
$ws.Range("J13").Value2 = 0.3333333333333333
$ws.Range("H15").Value2 = 0.1
$ws.Range("J15").Value2 = 0.7
$ws.Range("K15").Value2 = 0.03333333333333333
$ws.Range("O15").Value2 = 0.06666666666666667
$ws.Range("S15").Value2 = 0.1
$ws.Range("F16").Value2 = 0.08
$ws.Range("H16").Value2 = 0.04
$ws.Range("I16").Value2 = 0.04
$ws.Range("J16").Value2 = 0.4
$ws.Range("K16").Value2 = 0.04
$ws.Range("M16").Value2 = 0.04
$ws.Range("O16").Value2 = 0.12
$ws.Range("S16").Value2 = 0.24
$ws.Range("H17").Value2 = 0.1176470588235294
$ws.Range("I17").Value2 = 0.05882352941176471
$ws.Range("J17").Value2 = 0.611764705882353
$ws.Range("K17").Value2 = 0.07058823529411765
$ws.Range("M17").Value2 = 0.01176470588235294
$ws.Range("O17").Value2 = 0.05882352941176471
$ws.Range("S17").Value2 = 0.07058823529411765
$ws.Range("H18").Value2 = 0.2068965517241379
$ws.Range("I18").Value2 = 0.03448275862068965
$ws.Range("J18").Value2 = 0.5517241379310345
$ws.Range("K18").Value2 = 0.1206896551724138
$ws.Range("O18").Value2 = 0.01724137931034483
$ws.Range("S18").Value2 = 0.06896551724137931
$ws.Range("F19").Value2 = 0.02013422818791946
$ws.Range("H19").Value2 = 0.1140939597315436
$ws.Range("I19").Value2 = 0.02013422818791946
$ws.Range("J19").Value2 = 0.697986577181208
$ws.Range("K19").Value2 = 0.06040268456375839
$ws.Range("M19").Value2 = 0.006711409395973154
$ws.Range("O19").Value2 = 0.06040268456375839
$ws.Range("S19").Value2 = 0.02013422818791946
